$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("working")

# Update D6 value from 1000 to 100
$ws.Range("D6").Value = 100

# Add new row 7 with data
$ws.Range("A7").Value = 2022
$ws.Range("B7").Value = 1067
$ws.Range("C7").Value = 576
$ws.Range("D7").Value = 108.7

# Update selection to D8
$ws.Range("D8").Select()
